# Append two new paragraphs at the end of the document, right after the
# current last paragraph ("... Argentina"):
#   1. An empty paragraph
#   2. A paragraph containing "Samuelito el papelito"
# Both new paragraphs use the same run formatting (Helvetica, 24pt / sz=24)
# as the rest of the document.

$d = $word.ActiveDocument

$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

function New-PlainParagraphXml([string]$text) {
    # XML-escape the text content so arbitrary text is safe to embed.
    $escaped = $text.Replace('&', '&amp;').Replace('<', '&lt;').Replace('>', '&gt;')
    return '<w:p xmlns:w="' + $wNs + '"><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">' + $escaped + '</w:t></w:r></w:p>'
}

# Collapsed range positioned right at the end of the last paragraph's text
# (i.e. immediately before its paragraph mark) so InsertXML inserts new
# paragraphs there instead of overwriting anything.
$lastPara = $d.Paragraphs.Last
$insertPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$newParagraphsXml = (New-PlainParagraphXml '') + (New-PlainParagraphXml 'Samuelito el papelito')

$insertPoint.InsertXML($newParagraphsXml) | Out-Null
